$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new product row right before "STRINGAZOLE" (row 15), shifting rows 15-20 down by one.
$ws.Rows("15:15").Insert()

# Copy the formatting (styles, not values) from the row above so the new row matches
# the look of the other data rows in the table.
$ws.Range("A14:N14").Copy()
$ws.Range("A15:N15").PasteSpecial(-4122)
$ws.Rows("15:15").RowHeight = 25.5

# Re-create the merged cell groups for the new row (same layout as every other data row).
$ws.Range("B15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()

# Populate the new row with the new product's data.
$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "PANADOL COLD & FLU DAY 24 F.C. TABS"
$ws.Range("H15").Value = "2:1"
$ws.Range("L15").Value = 38
$ws.Range("N15").Value = "0:2"

# Renumber the serial numbers ("م" column) of the rows that shifted down.
$ws.Range("A16").Value = 13
$ws.Range("A17").Value = 14
$ws.Range("A18").Value = 15
$ws.Range("A19").Value = 16

# Restore the exact row heights used by the sheet (Excel's auto row-height after
# the insert/shift doesn't always keep each row's original value).
$ws.Rows("16:16").RowHeight = 24.75
$ws.Rows("17:17").RowHeight = 25.5
$ws.Rows("18:18").RowHeight = 25.5
$ws.Rows("19:19").RowHeight = 24.75

# Update the totals row: add the new product's amount to the running total.
$ws.Range("K20").Value = 314
$ws.Rows("20:20").RowHeight = 26.25
